# Auto-generated edit script: update market-price derived cells
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 113936.336
$ws.Range("I15").Value = 113936.336
$ws.Range("K15").Value = 341809.008
$ws.Range("M15").Value = -341640.008
$ws.Range("H40").Value = 2286
$ws.Range("I40").Value = 1866.6666
$ws.Range("K40").Value = 1866.6666
$ws.Range("M40").Value = -1691.6666
$ws.Range("H111").Value = 1157.1666
$ws.Range("I111").Value = 719.25
$ws.Range("J111").Value = 2033
$ws.Range("K111").Value = 2157.75
$ws.Range("L111").Value = 6099
$ws.Range("M111").Value = 909.25
$ws.Range("N111").Value = -12233
$ws.Range("H125").Value = 9343063
$ws.Range("I125").Value = 464
$ws.Range("J125").Value = 11211583
$ws.Range("K125").Value = 4176
$ws.Range("L125").Value = 100904247
$ws.Range("M125").Value = -1716
$ws.Range("N125").Value = -100909167
$ws.Range("H135").Value = 1134.8
$ws.Range("I135").Value = 1050.5807
$ws.Range("K135").Value = 9455.2263
$ws.Range("M135").Value = -6920.2263
$ws.Range("H137").Value = 1419.9706
$ws.Range("I137").Value = 1120.2
$ws.Range("J137").Value = 1471.6552
$ws.Range("K137").Value = 3360.6
$ws.Range("L137").Value = 4414.9656
$ws.Range("M137").Value = -810.6000000000004
$ws.Range("N137").Value = -9514.9656

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2708.6758
$ws.Range("I32").Value = 1919.5518
$ws.Range("J32").Value = 5569.25
$ws.Range("K32").Value = 1919.5518
$ws.Range("L32").Value = 5569.25
$ws.Range("M32").Value = -1632.5518
$ws.Range("N32").Value = -6143.25
$ws.Range("H45").Value = 1820.381
$ws.Range("I45").Value = 1275.1578
$ws.Range("K45").Value = 1275.1578
$ws.Range("M45").Value = -898.1578
$ws.Range("H74").Value = 1492.5518
$ws.Range("I74").Value = 1623.3334
$ws.Range("K74").Value = 1623.3334
$ws.Range("M74").Value = -749.3334
$ws.Range("H77").Value = 1492.5518
$ws.Range("I77").Value = 1623.3334
$ws.Range("K77").Value = 8116.666999999999
$ws.Range("M77").Value = -3748.666999999999
$ws.Range("H132").Value = 2040.4318
$ws.Range("I132").Value = 1637.2222
$ws.Range("J132").Value = 3854.875
$ws.Range("K132").Value = 4911.6666
$ws.Range("L132").Value = 11564.625
$ws.Range("M132").Value = -2381.6666
$ws.Range("N132").Value = -16624.625
$ws.Range("H133").Value = 53399.715
$ws.Range("J133").Value = 53399.715
$ws.Range("L133").Value = 53399.715
$ws.Range("N133").Value = -58459.715
$ws.Range("H139").Value = 62460.625
$ws.Range("J139").Value = 62460.625
$ws.Range("L139").Value = 62460.625
$ws.Range("N139").Value = -72740.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 46266.668
$ws.Range("J59").Value = 46266.668
$ws.Range("L59").Value = 46266.668
$ws.Range("N59").Value = -47960.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("H94").Value = 917.28125
$ws.Range("I94").Value = 990
$ws.Range("J94").Value = 900.5
$ws.Range("K94").Value = 990
$ws.Range("L94").Value = 900.5
$ws.Range("M94").Value = -539
$ws.Range("N94").Value = -1802.5
$ws.Range("H97").Value = 29999.5
$ws.Range("J97").Value = 29999.5
$ws.Range("L97").Value = 29999.5
$ws.Range("N97").Value = -31981.5
$ws.Range("H99").Value = 7813864
$ws.Range("I99").Value = 10417950
$ws.Range("J99").Value = 1607
$ws.Range("K99").Value = 10417950
$ws.Range("L99").Value = 1607
$ws.Range("M99").Value = -10416452
$ws.Range("N99").Value = -4603
$ws.Range("H109").Value = 42666.332
$ws.Range("J109").Value = 42666.332
$ws.Range("L109").Value = 42666.332
$ws.Range("N109").Value = -44746.332
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 747.375
$ws.Range("I122").Value = 658.38464
$ws.Range("J122").Value = 1133
$ws.Range("K122").Value = 1975.15392
$ws.Range("L122").Value = 3399
$ws.Range("M122").Value = 474.84608
$ws.Range("N122").Value = -8299
$ws.Range("H126").Value = 7813864
$ws.Range("I126").Value = 10417950
$ws.Range("J126").Value = 1607
$ws.Range("K126").Value = 31253850
$ws.Range("L126").Value = 4821
$ws.Range("M126").Value = -31251380
$ws.Range("N126").Value = -9761
$ws.Range("H132").Value = 2754.7585
$ws.Range("I132").Value = 1993.591
$ws.Range("K132").Value = 5980.772999999999
$ws.Range("M132").Value = -3450.772999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 15000
$ws.Range("N17").Value = -15338
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()
$ws.Range("N35").ClearContents()
$ws.Range("H68").Value = 1039.39
$ws.Range("I68").Value = 821.4808
$ws.Range("J68").Value = 1275.4584
$ws.Range("K68").Value = 2464.4424
$ws.Range("L68").Value = 3826.3752
$ws.Range("M68").Value = -1653.4424
$ws.Range("N68").Value = -5448.3752
$ws.Range("H71").Value = 1039.39
$ws.Range("I71").Value = 821.4808
$ws.Range("J71").Value = 1275.4584
$ws.Range("K71").Value = 7393.327200000001
$ws.Range("L71").Value = 11479.1256
$ws.Range("M71").Value = -3337.327200000001
$ws.Range("N71").Value = -19591.1256
$ws.Range("H107").Value = 1026.95
$ws.Range("I107").Value = 334.5
$ws.Range("J107").Value = 1323.7142
$ws.Range("K107").Value = 1003.5
$ws.Range("L107").Value = 3971.1426
$ws.Range("M107").Value = 916.5
$ws.Range("N107").Value = -7811.142599999999
$ws.Range("H133").Value = 11370.625
$ws.Range("I133").Value = 10393
$ws.Range("K133").Value = 31179
$ws.Range("M133").Value = -26119
$ws.Range("H137").Value = 4045314.2
$ws.Range("I137").Value = 9095532
$ws.Range("J137").Value = 77285.92999999999
$ws.Range("K137").Value = 27286596
$ws.Range("L137").Value = 231857.79
$ws.Range("M137").Value = -27281496
$ws.Range("N137").Value = -242057.79

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2092.8965
$ws.Range("I126").Value = 1679.4546
$ws.Range("J126").Value = 2345.5557
$ws.Range("K126").Value = 5038.3638
$ws.Range("L126").Value = 7036.6671
$ws.Range("M126").Value = -2568.3638
$ws.Range("N126").Value = -11976.6671
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 902.3570999999999
$ws.Range("I93").Value = 802.75
$ws.Range("K93").Value = 802.75
$ws.Range("M93").Value = 445.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 23543354
$ws.Range("I2").Value = 66675500
$ws.Range("J2").Value = 16728.092
$ws.Range("K2").Value = 66675500
$ws.Range("L2").Value = 16728.092
$ws.Range("M2").Value = -66675388
$ws.Range("N2").Value = -16952.092
$ws.Range("H4").Value = 16870.375
$ws.Range("I4").Value = 2653.3333
$ws.Range("J4").Value = 25400.6
$ws.Range("K4").Value = 2653.3333
$ws.Range("L4").Value = 25400.6
$ws.Range("M4").Value = -2540.3333
$ws.Range("N4").Value = -25626.6
$ws.Range("H47").Value = 13752.286
$ws.Range("J47").Value = 13752.286
$ws.Range("L47").Value = 13752.286
$ws.Range("N47").Value = -14896.286
$ws.Range("H86").Value = 14375
$ws.Range("J86").Value = 14375
$ws.Range("L86").Value = 14375
$ws.Range("N86").Value = -16621
$ws.Range("H89").Value = 14375
$ws.Range("J89").Value = 14375
$ws.Range("L89").Value = 71875
$ws.Range("N89").Value = -83107
$ws.Range("H107").Value = 679.1111
$ws.Range("I107").Value = 663.3125
$ws.Range("J107").Value = 805.5
$ws.Range("K107").Value = 1989.9375
$ws.Range("L107").Value = 2416.5
$ws.Range("M107").Value = -69.9375
$ws.Range("N107").Value = -6256.5
$ws.Range("H126").Value = 126273.875
$ws.Range("I126").Value = 144027.28
$ws.Range("K126").Value = 432081.84
$ws.Range("M126").Value = -429611.84
